$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing data range
$ws.Range("A1:D2").ClearContents()

# Header row
$ws.Range("A1").Value = "a"
$ws.Range("B1").Value = "b"

# Row 2: only A
$ws.Range("A2").Value = "first of a"

# Row 3: only B
$ws.Range("B3").Value = "first of b"

# Row 4: only A
$ws.Range("A4").Value = "second of only a"

# Row 5: both A and B
$ws.Range("A5").Value = "a of has both"
$ws.Range("B5").Value = "b of has both"
